$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 0.1425304360315636
    "D2" = 0.1058422808887993
    "E2" = 0.13028179933346
    "F2" = 2.878478165862106
    "G2" = 2.368769996223534
    "H2" = 1.696868061313182
    "J2" = 0.2151472402880472
    "N2" = 1.995241137945641
    "B3" = 0.1329091933101267
    "D3" = 0.1057099897463161
    "E3" = 0.129376248143025
    "F3" = 2.750918434283136
    "G3" = 2.222640816841505
    "H3" = 1.636692919714051
    "J3" = 0.208773148833032
    "N3" = 1.846526605564776
    "B4" = 0.1270661492633423
    "D4" = 0.1056316089112226
    "E4" = 0.1288919000434596
    "F4" = 2.674709327376632
    "G4" = 2.13468666062019
    "H4" = 1.600961959854857
    "J4" = 0.205012059613054
    "N4" = 1.755520437593759
    "B5" = 0.1247015078377984
    "D5" = 0.1056003683078792
    "E5" = 0.1287124574646867
    "F5" = 2.64417696227224
    "G5" = 2.099280920985308
    "H5" = 1.586703403558658
    "J5" = 0.2035173224339673
    "N5" = 1.71851688491077
    "B6" = 0.1243098614099978
    "D6" = 0.1055952226743999
    "E6" = 0.128683741779934
    "F6" = 2.639138474748648
    "G6" = 2.093427924020915
    "H6" = 1.584353915396832
    "J6" = 0.2032714010193857
    "N6" = 1.712377599108521
    "B7" = 0.1270341919676241
    "D7" = 0.1056311847747633
    "E7" = 0.1288894075201767
    "F7" = 2.674295448035082
    "G7" = 2.134207411306051
    "H7" = 1.600768445849013
    "J7" = 0.2049917480354679
    "N7" = 1.755021054646477
    "B8" = 0.1391997998659633
    "D8" = 0.1057960665953246
    "E8" = 0.129954637554043
    "F8" = 2.834052162487836
    "G8" = 2.318011996819791
    "H8" = 1.675864740954125
    "J8" = 0.212917504022812
    "N8" = 1.94390425161825
    "B9" = 0.1635584291339143
    "D9" = 0.1061427228871796
    "E9" = 0.1326164825695173
    "F9" = 3.164481609023056
    "G9" = 2.692914156986546
    "H9" = 1.832970800449971
    "J9" = 0.2296922894002904
    "N9" = 2.316504389646809
    "B10" = 0.1817500012197826
    "D10" = 0.106412707435446
    "E10" = 0.1349279334950673
    "F10" = 3.418265018242977
    "G10" = 2.977790851650184
    "H10" = 1.954674957977772
    "J10" = 0.2427996710749767
    "N10" = 2.591327722520646
    "B11" = 0.1900877858894034
    "D11" = 0.106539096441006
    "E11" = 0.1360582219594733
    "F11" = 3.536236937461979
    "G11" = 3.10957700108986
    "H11" = 2.011467405619328
    "J11" = 0.2489397106595561
    "N11" = 2.716530669508472
    "B12" = 0.1932538308942782
    "D12" = 0.1065874905544781
    "E12" = 0.1364976832432454
    "F12" = 3.581282987014646
    "G12" = 3.15980762891769
    "H12" = 2.033183669984112
    "J12" = 0.2512908683668513
    "N12" = 2.76396319373481
    "B13" = 0.1925715833274921
    "D13" = 0.1065770439752107
    "E13" = 0.1364025264829394
    "F13" = 3.571564784942183
    "G13" = 3.148974891615353
    "H13" = 2.02849725803037
    "J13" = 0.2507833367248367
    "N13" = 2.753746911386088
    "B14" = 0.1903480852990072
    "D14" = 0.1065430670494472
    "E14" = 0.1360941466472845
    "F14" = 3.539935378736061
    "G14" = 3.113702908760729
    "H14" = 2.013249775342672
    "J14" = 0.2491326164057455
    "N14" = 2.720432585385765
    "B15" = 0.1889872547793345
    "D15" = 0.1065223252554794
    "E15" = 0.1359067491117649
    "F15" = 3.520610273165971
    "G15" = 3.092140582485115
    "H15" = 2.003937769423828
    "J15" = 0.2481249124802929
    "N15" = 2.700029152110062
    "B16" = 0.1812063398621291
    "D16" = 0.1064045211602469
    "E16" = 0.1348556618023871
    "F16" = 3.410606892642278
    "G16" = 2.969223366258007
    "H16" = 1.950992589575719
    "J16" = 0.2424020215006948
    "N16" = 2.583148644060884
    "B17" = 0.1764487851816483
    "D17" = 0.1063331807884911
    "E17" = 0.1342311160739662
    "F17" = 3.343776744307405
    "G17" = 2.894387060466613
    "H17" = 1.918881687400301
    "J17" = 0.2389370388622893
    "N17" = 2.511489453097283
    "B18" = 0.1737182494253204
    "D18" = 0.1062924826548013
    "E18" = 0.1338793036318329
    "F18" = 3.305574844941731
    "G18" = 2.851549089412629
    "H18" = 1.900546521073124
    "J18" = 0.2369607392022033
    "N18" = 2.470290675909951
    "B19" = 0.1727947553761737
    "D19" = 0.1062787599065302
    "E19" = 0.1337614551878872
    "F19" = 3.292680760349469
    "G19" = 2.837079934219787
    "H19" = 1.89436145048802
    "J19" = 0.2362944442791957
    "N19" = 2.456344689235948
    "B20" = 0.1769546279148955
    "D20" = 0.1063407402842653
    "E20" = 0.1342968323355507
    "F20" = 3.350866326720222
    "G20" = 2.90233212842918
    "H20" = 1.922286027038183
    "J20" = 0.2393041633429505
    "N20" = 2.519115905984222
    "B21" = 0.1910009463252038
    "D21" = 0.1065530322609014
    "E21" = 0.136184413718567
    "F21" = 3.549215510681506
    "G21" = 3.124054210283646
    "H21" = 2.017722587296021
    "J21" = 0.249616761358638
    "N21" = 2.730217288371648
    "B22" = 0.2002316580428385
    "D22" = 0.1066948957068323
    "E22" = 0.1374848268051494
    "F22" = 3.681024872415463
    "G22" = 3.270868098814447
    "H22" = 2.081323391753187
    "J22" = 0.2565087896800833
    "N22" = 2.868303394937357
    "B23" = 0.195300506148314
    "D23" = 0.1066188884862704
    "E23" = 0.1367846235793699
    "F23" = 3.610473322205905
    "G23" = 3.192332840114375
    "H23" = 2.047264520717476
    "J23" = 0.2528162789516131
    "N23" = 2.794595190493112
    "B24" = 0.1767259219250263
    "D24" = 0.10633732165204
    "E24" = 0.1342670994488167
    "F24" = 3.347660443095805
    "G24" = 2.898739583346128
    "H24" = 1.920746533548197
    "J24" = 0.2391381373017794
    "N24" = 2.515667989176507
    "B25" = 0.1569161479393415
    "D25" = 0.1060463370010964
    "E25" = 0.1318343853988502
    "F25" = 3.073198553213103
    "G25" = 2.589881758271019
    "H25" = 1.789387422476921
    "J25" = 0.2250191059075775
    "N25" = 2.215498394811959
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
